$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 59,5
$data[0,0] = 56
$data[0,1] = 847
$data[0,2] = 17
$data[0,3] = 77
$data[0,4] = 86
$data[1,0] = 53
$data[1,1] = 847
$data[1,2] = 131
$data[1,3] = 79
$data[1,4] = 107
$data[2,0] = 46
$data[2,1] = 852
$data[2,2] = 259
$data[2,3] = 61
$data[2,4] = 104
$data[3,0] = 59
$data[3,1] = 682
$data[3,2] = 0
$data[3,3] = 88
$data[3,4] = 119
$data[4,0] = 50
$data[4,1] = 682
$data[4,2] = 140
$data[4,3] = 80
$data[4,4] = 100
$data[5,0] = 44
$data[5,1] = 676
$data[5,2] = 263
$data[5,3] = 90
$data[5,4] = 96
$data[6,0] = 41
$data[6,1] = 685
$data[6,2] = 378
$data[6,3] = 74
$data[6,4] = 116
$data[7,0] = 34
$data[7,1] = 685
$data[7,2] = 517
$data[7,3] = 59
$data[7,4] = 92
$data[8,0] = 29
$data[8,1] = 685
$data[8,2] = 643
$data[8,3] = 74
$data[8,4] = 78
$data[9,0] = 25
$data[9,1] = 685
$data[9,2] = 756
$data[9,3] = 76
$data[9,4] = 100
$data[10,0] = 20
$data[10,1] = 679
$data[10,2] = 882
$data[10,3] = 82
$data[10,4] = 111
$data[11,0] = 15
$data[11,1] = 684
$data[11,2] = 1022
$data[11,3] = 75
$data[11,4] = 83
$data[12,0] = 10
$data[12,1] = 685
$data[12,2] = 1140
$data[12,3] = 70
$data[12,4] = 90
$data[13,0] = 5
$data[13,1] = 687
$data[13,2] = 1265
$data[13,3] = 61
$data[13,4] = 97
$data[14,0] = 1
$data[14,1] = 686
$data[14,2] = 1393
$data[14,3] = 73
$data[14,4] = 84
$data[15,0] = 57
$data[15,1] = 523
$data[15,2] = 10
$data[15,3] = 72
$data[15,4] = 111
$data[16,0] = 49
$data[16,1] = 527
$data[16,2] = 142
$data[16,3] = 70
$data[16,4] = 93
$data[17,0] = 47
$data[17,1] = 516
$data[17,2] = 254
$data[17,3] = 81
$data[17,4] = 110
$data[18,0] = 38
$data[18,1] = 523
$data[18,2] = 393
$data[18,3] = 76
$data[18,4] = 87
$data[19,0] = 33
$data[19,1] = 520
$data[19,2] = 520
$data[19,3] = 84
$data[19,4] = 79
$data[20,0] = 32
$data[20,1] = 535
$data[20,2] = 620
$data[20,3] = 42
$data[20,4] = 67
$data[21,0] = 27
$data[21,1] = 523
$data[21,2] = 691
$data[21,3] = 58
$data[21,4] = 58
$data[22,0] = 23
$data[22,1] = 519
$data[22,2] = 766
$data[22,3] = 66
$data[22,4] = 95
$data[23,0] = 14
$data[23,1] = 516
$data[23,2] = 1116
$data[23,3] = 76
$data[23,4] = 119
$data[24,0] = 8
$data[24,1] = 518
$data[24,2] = 1250
$data[24,3] = 59
$data[24,4] = 104
$data[25,0] = 3
$data[25,1] = 520
$data[25,2] = 1376
$data[25,3] = 46
$data[25,4] = 93
$data[26,0] = 58
$data[26,1] = 368
$data[26,2] = 10
$data[26,3] = 78
$data[26,4] = 116
$data[27,0] = 51
$data[27,1] = 355
$data[27,2] = 137
$data[27,3] = 89
$data[27,4] = 100
$data[28,0] = 43
$data[28,1] = 361
$data[28,2] = 265
$data[28,3] = 76
$data[28,4] = 103
$data[29,0] = 40
$data[29,1] = 370
$data[29,2] = 389
$data[29,3] = 65
$data[29,4] = 81
$data[30,0] = 37
$data[30,1] = 364
$data[30,2] = 502
$data[30,3] = 67
$data[30,4] = 94
$data[31,0] = 31
$data[31,1] = 366
$data[31,2] = 625
$data[31,3] = 61
$data[31,4] = 108
$data[32,0] = 26
$data[32,1] = 360
$data[32,2] = 752
$data[32,3] = 80
$data[32,4] = 113
$data[33,0] = 19
$data[33,1] = 359
$data[33,2] = 890
$data[33,3] = 79
$data[33,4] = 98
$data[34,0] = 17
$data[34,1] = 362
$data[34,2] = 1006
$data[34,3] = 59
$data[34,4] = 99
$data[35,0] = 11
$data[35,1] = 358
$data[35,2] = 1128
$data[35,3] = 66
$data[35,4] = 94
$data[36,0] = 9
$data[36,1] = 357
$data[36,2] = 1240
$data[36,3] = 80
$data[36,4] = 113
$data[37,0] = 4
$data[37,1] = 355
$data[37,2] = 1362
$data[37,3] = 66
$data[37,4] = 107
$data[38,0] = 54
$data[38,1] = 190
$data[38,2] = 21
$data[38,3] = 80
$data[38,4] = 86
$data[39,0] = 48
$data[39,1] = 193
$data[39,2] = 145
$data[39,3] = 75
$data[39,4] = 90
$data[40,0] = 42
$data[40,1] = 190
$data[40,2] = 377
$data[40,3] = 80
$data[40,4] = 108
$data[41,0] = 35
$data[41,1] = 189
$data[41,2] = 512
$data[41,3] = 75
$data[41,4] = 87
$data[42,0] = 28
$data[42,1] = 194
$data[42,2] = 655
$data[42,3] = 70
$data[42,4] = 73
$data[43,0] = 22
$data[43,1] = 189
$data[43,2] = 879
$data[43,3] = 73
$data[43,4] = 100
$data[44,0] = 16
$data[44,1] = 198
$data[44,2] = 1013
$data[44,3] = 65
$data[44,4] = 83
$data[45,0] = 12
$data[45,1] = 205
$data[45,2] = 1128
$data[45,3] = 39
$data[45,4] = 95
$data[46,0] = 7
$data[46,1] = 190
$data[46,2] = 1255
$data[46,3] = 76
$data[46,4] = 90
$data[47,0] = 2
$data[47,1] = 193
$data[47,2] = 1380
$data[47,3] = 74
$data[47,4] = 78
$data[48,0] = 55
$data[48,1] = 31
$data[48,2] = 19
$data[48,3] = 82
$data[48,4] = 83
$data[49,0] = 52
$data[49,1] = 27
$data[49,2] = 136
$data[49,3] = 88
$data[49,4] = 99
$data[50,0] = 45
$data[50,1] = 27
$data[50,2] = 261
$data[50,3] = 87
$data[50,4] = 110
$data[51,0] = 39
$data[51,1] = 34
$data[51,2] = 392
$data[51,3] = 87
$data[51,4] = 86
$data[52,0] = 36
$data[52,1] = 41
$data[52,2] = 509
$data[52,3] = 42
$data[52,4] = 92
$data[53,0] = 30
$data[53,1] = 35
$data[53,2] = 633
$data[53,3] = 73
$data[53,4] = 96
$data[54,0] = 24
$data[54,1] = 33
$data[54,2] = 762
$data[54,3] = 75
$data[54,4] = 84
$data[55,0] = 21
$data[55,1] = 29
$data[55,2] = 882
$data[55,3] = 79
$data[55,4] = 84
$data[56,0] = 18
$data[56,1] = 22
$data[56,2] = 1001
$data[56,3] = 85
$data[56,4] = 92
$data[57,0] = 13
$data[57,1] = 25
$data[57,2] = 1120
$data[57,3] = 81
$data[57,4] = 117
$data[58,0] = 6
$data[58,1] = 32
$data[58,2] = 1259
$data[58,3] = 88
$data[58,4] = 82

$rng = $ws.Range("A2:E60")
$rng.Value = $data
